# Applies the template-guideline.docx style update:
#   - Title style: paragraph alignment changes from Left to Center.
#   - Subtitle style: keeps (effectively inherits/matches) a Center
#     alignment now that its base style (Title) is itself Center.
#
# wdAlignParagraphCenter = 1 (WdParagraphAlignment)
$wdAlignParagraphCenter = 1

$d = $word.ActiveDocument

# --- Title paragraph style: left -> center -------------------------------
$titleStyle = $d.Styles("Title")
$titleStyle.ParagraphFormat.Alignment = $wdAlignParagraphCenter

# --- Subtitle paragraph style: stays centered -----------------------------
# Subtitle is basedOn Title, so it now matches its base style's alignment.
$subtitleStyle = $d.Styles("Subtitle")
$subtitleStyle.ParagraphFormat.Alignment = $wdAlignParagraphCenter
